$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: bump the "Last status check on" timestamp (15:15 -> 15:30)
$ws.Range("F1").Value = "Last status check on: 14.01.2022 15:30"

# Row 2 (TankONO) got a fresh scrape: current price rose to 33.9,
# the previous current price (33.5) slides into the "Old Cena" column,
# the delta is now rendered as literal text "+0.4" and the old
# numeric-date timestamp is replaced by a plain text timestamp string.
$ws.Range("B2").Value = 33.9
$ws.Range("C2").Value = 33.5

# Force D2 to hold literal text "+0.4" (a leading apostrophe stops Excel
# from re-parsing it back into the number 0.4).
$ws.Range("D2").Value = "'+0.4"

# E2 becomes a plain text timestamp (no longer a numeric date serial),
# so reset it to the default "Normal" style before writing the string.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "2022-01-14 15:30:04"
